$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The import's last data row (row 18) is moved up to become the first data
# row (row 2), pushing the existing rows 2-17 down to 3-18 - i.e. row 18 is
# promoted to the top of the data block.
$lastRow = 18
$colCount = 9

# 1) Capture the raw (non-date-coerced) values of the row being moved before
#    any shifting happens.
$values = @{}
for ($c = 1; $c -le $colCount; $c++) {
    $values[$c] = $ws.Cells.Item($lastRow, $c).Value2()
}

# 2) Insert a blank row at row 2; this shifts rows 2-18 down to rows 3-19
#    (the row we captured above is now sitting at $lastRow + 1).
$ws.Rows.Item(2).Insert(-4121) | Out-Null

# 3) Copy the formatting (styles) of the shifted source row onto the new
#    row 2, restricted to the used columns (A:I) so we don't touch the
#    rest of the (now effectively infinite) row.
$srcRange = $ws.Range($ws.Cells.Item($lastRow + 1, 1), $ws.Cells.Item($lastRow + 1, $colCount))
$dstRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $colCount))
$srcRange.Copy() | Out-Null
$dstRange.PasteSpecial(-4122) | Out-Null

# 4) Write the captured values into the new row 2 (clearing any column that
#    was genuinely empty on the source row, e.g. column H).
for ($c = 1; $c -le $colCount; $c++) {
    $v = $values[$c]
    if ($v -ne $null -and $v -ne "") {
        $ws.Cells.Item(2, $c).Value2 = $v
    } else {
        $ws.Cells.Item(2, $c).Clear() | Out-Null
    }
}

# 5) Remove the now-duplicated old row (original row 18, shifted to
#    $lastRow + 1) since its contents now live at row 2.
$ws.Rows.Item($lastRow + 1).Delete() | Out-Null

# 6) Match the resulting selection: the whole of (now relocated) row 2.
$ws.Range("A2:XFD2").Select() | Out-Null
